$wb = $excel.ActiveWorkbook
$originalActiveSheet = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the "总计" (total) sheet, and
#    populate it with the per-fund holdings for that quarter (same layout as
#    the other quarterly sheets: 2020-Q4 / 2021-Q1 / 2021-Q3 / 2021-Q4).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(4)      # "2021-Q4" - a sheet with the fund-table layout/styles
$totalSheet = $wb.Worksheets.Item(5)    # "总计" - new sheet goes right before this one

$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Reuse the header-row (B1:H1) formatting and the column-A index-number
# formatting (bold/bordered style) from the template sheet.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "040018"
$q1.Range("C2").Value = "华安香港精选股票(QDII)"
$q1.Range("D2").Value = "5.47"
$q1.Range("E2").Value = "88.46"
$q1.Range("F2").Value = "2.46"
$q1.Range("G2").Value = "0.1346"
$q1.Range("H2").Value = 9

$q1.Range("A3").Value = 1
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = "040021"
$q1.Range("C3").Value = "华安大中华升级股票(QDII)"
$q1.Range("D3").Value = "0.26"
$q1.Range("E3").Value = "87.37"
$q1.Range("F3").Value = "2.34"
$q1.Range("G3").Value = "0.0061"
$q1.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# 2) Add a "2022-Q1" row to the top of the "总计" (total) summary sheet,
#    pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$b2 = $total.Range("B2").Value()
$c2 = $total.Range("C2").Value()
$d2 = $total.Range("D2").Value()
$b3 = $total.Range("B3").Value()
$c3 = $total.Range("C3").Value()
$d3 = $total.Range("D3").Value()
$b4 = $total.Range("B4").Value()
$c4 = $total.Range("C4").Value()
$d4 = $total.Range("D4").Value()
$b5 = $total.Range("B5").Value()
$c5 = $total.Range("C5").Value()
$d5 = $total.Range("D5").Value()

# extend the column-A index-number styling down to the new last row
$total.Range("A2").Copy()
$total.Range("A6").PasteSpecial(-4122)

$total.Range("A6").Value = 4
$total.Range("B6").Value = $b5
$total.Range("C6").Value = $c5
$total.Range("D6").Value = $d5

$total.Range("A5").Value = 3
$total.Range("B5").Value = $b4
$total.Range("C5").Value = $c4
$total.Range("D5").Value = $d4

$total.Range("A4").Value = 2
$total.Range("B4").Value = $b3
$total.Range("C4").Value = $c3
$total.Range("D4").Value = $d3

$total.Range("A3").Value = 1
$total.Range("B3").Value = $b2
$total.Range("C3").Value = $c2
$total.Range("D3").Value = $d2

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.14

# Restore the originally-active sheet/selection so this edit doesn't churn
# unrelated UI-state (active tab) beyond what the diff describes.
$originalActiveSheet.Activate()

